$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: a new price record is prepended at row 100 and every
# existing record from row 100 down shifts one row lower (the oldest record,
# previously at row 213, ends up at the new row 214).
$ws.Rows.Item(100).Insert()

# Populate the new row 100 with the latest weekly price record.
$ws.Cells.Item(100, 1).Value = 10
$ws.Cells.Item(100, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(100, 3).Value = 'La Araucanía'
$ws.Cells.Item(100, 4).Value = 44893
$ws.Cells.Item(100, 5).Value = 9
$ws.Cells.Item(100, 6).Value = 100112012
$ws.Cells.Item(100, 7).Value = 'Espinaca'
$ws.Cells.Item(100, 8).Value = 'Sin especificar'
$ws.Cells.Item(100, 9).Value = 'Primera'
$ws.Cells.Item(100, 10).Value = 40
$ws.Cells.Item(100, 11).Value = 9000
$ws.Cells.Item(100, 12).Value = 9000
$ws.Cells.Item(100, 13).Value = 9000
$ws.Cells.Item(100, 14).Value = '$/docena de atados'
$ws.Cells.Item(100, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(100, 16).Value = 3000
$ws.Cells.Item(100, 17).Value = 3
$ws.Cells.Item(100, 18).Value = 'Hortaliza'
